# Apply the LOQ4071.xlsx content update:
#  - Row 10 (Objetivos:) gets a new "Objectives" body paragraph in B/C.
#  - A new row is inserted after row 12 (Docentes responsaveis:) so that the
#    "5840535 - Messias Borges Silva" answer moves to its own row (13),
#    "Programa resumido:" moves to row 14 with a new body text, and
#    everything that used to be on rows 13-23 shifts down to 14-24, with a
#    couple of new/changed text bodies along the way.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$xlPasteFormats = -4122

# --- 1. Update the "Objetivos:" body text (row 10, columns B & C) ---------
$ws.Cells.Item(10, 2).Value = "Complementar a formação multidisciplinar dos alunos de Engenharia abordando, com maior profundidade, tópicos atuais e relevantes sobre gestão da qualidade."
$ws.Cells.Item(10, 3).Value = "Complementar a formação multidisciplinar dos alunos de Engenharia abordando, com maior profundidade, tópicos atuais e relevantes sobre gestão da qualidade."

# --- 2. Insert a blank row before the old row 13 ---------------------------
# This shifts old rows 13..23 down to 14..24 (carrying their row heights and
# formatting along for free).
$ws.Rows.Item(13).Insert()

# The insert can leave a phantom styled-but-empty cell in A13; remove it so
# that row 13 only contains B13/C13, as in the target layout.
$ws.Cells.Item(13, 1).Clear()

# Newly-materialised B13/C13 cells would otherwise inherit a stray "bold"
# look from the row above; explicitly copy the correct (non-bold / red)
# column formatting down from row 14 before putting values in them.
$ws.Cells.Item(14, 2).Copy()
$ws.Cells.Item(13, 2).PasteSpecial($xlPasteFormats)
$ws.Cells.Item(14, 3).Copy()
$ws.Cells.Item(13, 3).PasteSpecial($xlPasteFormats)

# --- 3. Row 13: the "5840535 - Messias Borges Silva" line (B & C only) ----
$ws.Cells.Item(13, 2).Value = "5840535 - Messias Borges Silva"
$ws.Cells.Item(13, 3).Value = "5840535 - Messias Borges Silva"

# --- 4. Row 14: "Programa resumido:" + new body text -----------------------
$ws.Cells.Item(14, 1).Value = "Programa resumido:"
$ws.Cells.Item(14, 2).Value = "A definir, de acordo com o tópico programado."
$ws.Cells.Item(14, 3).Value = "A definir, de acordo com o tópico programado."

# --- 5. Row 15: "Short syllabus:" only -------------------------------------
$ws.Cells.Item(15, 1).Value = "Short syllabus:"
$ws.Cells.Item(15, 2).Clear()
$ws.Cells.Item(15, 3).Clear()

# --- 6. Row 16: "Programa:" + new body text --------------------------------
$ws.Cells.Item(16, 1).Value = "Programa:"
$ws.Cells.Item(16, 2).Value = "O conteúdo desta disciplina será de acordo com o tópico a ser programado, devendo abordar assuntos complementares a formação de um profissional de Engenharia."
$ws.Cells.Item(16, 3).Value = "O conteúdo desta disciplina será de acordo com o tópico a ser programado, devendo abordar assuntos complementares a formação de um profissional de Engenharia."

# --- 7. Row 17: "Syllabus:" only -------------------------------------------
$ws.Cells.Item(17, 1).Value = "Syllabus:"

# --- 8. Row 18: "Avaliação:" only, body text removed -----------------------
$ws.Cells.Item(18, 1).Value = "Avaliação:"
$ws.Cells.Item(18, 2).Clear()
$ws.Cells.Item(18, 3).Clear()

# --- 9. Row 19: "Método:" + its existing body text --------------------------
$ws.Cells.Item(19, 1).Value = "Método:"
$ws.Cells.Item(19, 2).Value = "O desenvolvimento da disciplina será baseado em leituras, aula expositiva, discussão e resolução de estudos de caso e resolução de exercícios."
$ws.Cells.Item(19, 3).Value = "O desenvolvimento da disciplina será baseado em leituras, aula expositiva, discussão e resolução de estudos de caso e resolução de exercícios."

# --- 10. Row 20: "Critério:" + "Provas e trabalhos." -----------------------
$ws.Cells.Item(20, 1).Value = "Critério:"
$ws.Cells.Item(20, 2).Value = "Provas e trabalhos."
$ws.Cells.Item(20, 3).Value = "Provas e trabalhos."

# --- 11. Row 21: "Norma de recuperação:" + recovery rule -------------------
$ws.Cells.Item(21, 1).Value = "Norma de recuperação:"
$ws.Cells.Item(21, 2).Value = "Prova única com nota maior ou igual a 5,0 (cinco)."
$ws.Cells.Item(21, 3).Value = "Prova única com nota maior ou igual a 5,0 (cinco)."

# --- 12. Row 22: "Bibliografia:" + new bibliography text -------------------
$ws.Cells.Item(22, 1).Value = "Bibliografia:"
$ws.Cells.Item(22, 2).Value = "Textos fornecidos pelo professor da disciplina`nArtigos extraídos de revistas especializadas na área de gestão e produção."
$ws.Cells.Item(22, 3).Value = "Textos fornecidos pelo professor da disciplina`nArtigos extraídos de revistas especializadas na área de gestão e produção."

# --- 13. Row 23: "Requisitos:" only -----------------------------------------
$ws.Cells.Item(23, 1).Value = "Requisitos:"
$ws.Cells.Item(23, 2).Clear()
$ws.Cells.Item(23, 3).Clear()

# --- 14. Row 24: requisite detail text (B & C only) -------------------------
$ws.Cells.Item(24, 1).Clear()
$ws.Cells.Item(24, 2).Value = "LOQ4044 -  Introdução à Engenharia da Qualidade  (Requisito fraco)`n"
$ws.Cells.Item(24, 3).Value = "LOQ4044 -  Introdução à Engenharia da Qualidade  (Requisito fraco)`n"

Write-Output "done"
